$d = $word.ActiveDocument

$replacements = @(
    @("47÷9=5, 2", "19÷8=2, 3"),
    @("14÷5=2, 4", "25÷3=8, 1"),
    @("23÷4=5, 3", "96÷9=10, 6"),
    @("13÷2=6, 1", "83÷6=13, 5"),
    @("64÷4=16, 0", "68÷6=11, 2"),
    @("43÷5=8, 3", "12÷9=1, 3"),
    @("76÷6=12, 4", "30÷7=4, 2"),
    @("44÷5=8, 4", "43÷5=8, 3"),
    @("82÷9=9, 1", "95÷9=10, 5"),
    @("77÷4=19, 1", "59÷5=11, 4"),
    @("12÷5=2, 2", "31÷4=7, 3"),
    @("32÷7=4, 4", "19÷2=9, 1"),
    @("66÷2=33, 0", "43÷2=21, 1"),
    @("88÷9=9, 7", "54÷8=6, 6"),
    @("99÷7=14, 1", "69÷6=11, 3"),
    @("30÷8=3, 6", "11÷8=1, 3"),
    @("27÷4=6, 3", "63÷9=7, 0"),
    @("37÷6=6, 1", "32÷7=4, 4"),
    @("60÷5=12, 0", "12÷6=2, 0"),
    @("66÷8=8, 2", "96÷3=32, 0"),
    @("18÷6=3, 0", "27÷9=3, 0"),
    @("58÷9=6, 4", "41÷8=5, 1"),
    @("52÷2=26, 0", "98÷8=12, 2"),
    @("58÷2=29, 0", "14÷9=1, 5"),
    @("47÷7=6, 5", "30÷5=6, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
